# Add four new "sign-out" / "unknown login" related translation rows
# to the "Import" sheet (rows 38-41), mirroring the existing cs (Czech)
# translation rows above them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Import")

# Use the last existing data row (37) as a formatting template so the new
# rows inherit the same cell style ("import": wrap text etc.) as the rest
# of the table.
$templateRange = $ws.Range($ws.Cells.Item(37, 1), $ws.Cells.Item(37, 3))
for ($r = 38; $r -le 41; $r++) {
    $targetRange = $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 3))
    $templateRange.Copy($targetRange)
}

# Column A ("language") is the same "cs" value used throughout the sheet.
$ws.Cells.Item(38, 1).Value = "cs"
$ws.Cells.Item(39, 1).Value = "cs"
$ws.Cells.Item(40, 1).Value = "cs"
$ws.Cells.Item(41, 1).Value = "cs"

# Keys (column B) for the sign-out screen.
$ws.Cells.Item(38, 2).Value = "public.sign-out.title"
$ws.Cells.Item(39, 2).Value = "public.sign-out"

# Translations (column C) for the sign-out screen.
$ws.Cells.Item(38, 3).Value = "Odhlašování"
$ws.Cells.Item(39, 3).Value = "Probíhá odhlašování z aplikace, prosím vyčkejte…"

# Key + translation for the "unknown user" error.
$ws.Cells.Item(40, 2).Value = "error.Who are you?"
$ws.Cells.Item(40, 3).Value = "Je nám líto, ale aplikace vás nepoznává."

# Key + translation for the "unknown login" error.
$ws.Cells.Item(41, 2).Value = "error.Unknown login"
$ws.Cells.Item(41, 3).Value = "Přihlášení selhalo, zkontrolujte si prosím jméno a heslo."

# Move the cursor/selection and scroll position to match the author's
# final view of the sheet.
$ws.Activate()
[void]$ws.Cells.Item(34, 2).Select()
$excel.ActiveWindow.ScrollRow = 19
